$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Stage 1")

# Insert a blank row above the footer/header row (row 59), shifting the
# header (row 59->60) and the legend row (row 60->61) down by one.
$ws.Range("A59:AJ59").Insert(-4121)  # xlShiftDown

# The new row continues the alternating color-band formatting of the data
# rows above it, matching row 47's format (same band as the new row).
$ws.Range("A47:AJ47").Copy()
$ws.Range("A59:AJ59").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row's values ("T55" state).
$ws.Range("A59").Value = "T55"
$vals = @(1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,1,1,0,0,0,0,1,1,1,0,1,1,0,0,0,0,0,0,0)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "59").Value = $vals[$i]
}

# Restore the view: scrolled near the bottom of the table with R60 selected.
$ws.Application.ActiveWindow.ScrollRow = 46
$ws.Range("R60").Select()
